$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'34.081.68"
$ws.Range("E2").Value = "  +0.11%  "

# Row 3
$ws.Range("D3").Value = "'1.788.69"
$ws.Range("E3").Value = "  +0.38%  "

# Row 5
$ws.Range("D5").Value = "'227.32"
$ws.Range("E5").Value = "  +1.45%  "

# Row 6
$ws.Range("E6").Value = "  -0.34%  "

# Row 8
$ws.Range("D8").Value = "'32.27"
$ws.Range("E8").Value = "  -0.39%  "

# Row 9
$ws.Range("E9").Value = "  +3.88%  "

# Row 10
$ws.Range("D10").Value = "'0.0687"
$ws.Range("E10").Value = "  -2.11%  "

# Row 11
$ws.Range("E11").Value = "  +0.88%  "

# Row 12
$ws.Range("E12").Value = "  +0.32%  "

# Row 13
$ws.Range("E13").Value = "  +5.52%  "

# Row 14
$ws.Range("D14").Value = "'1.800.41"
$ws.Range("E14").Value = "  +1.09%  "

# Row 15
$ws.Range("D15").Value = "'0.622"
$ws.Range("E15").Value = "  +0.32%  "

# Row 16
$ws.Range("D16").Value = "'34.050.28"
$ws.Range("E16").Value = "  +0.11%  "

# Row 17
$ws.Range("D17").Value = "'4.18"
$ws.Range("E17").Value = "  +1.08%  "

# Row 18
$ws.Range("D18").Value = "'68.12"
$ws.Range("E18").Value = "  +0.83%  "

# Row 19
$ws.Range("D19").Value = "'243.35"
$ws.Range("E19").Value = "  +0.57%  "

# Row 20
$ws.Range("E20").Value = "  +0.15%  "

# Row 21
$ws.Range("E21").Value = "  +0.04%  "

# Row 22
$ws.Range("D22").Value = "'10.87"
$ws.Range("E22").Value = "  +2.04%  "

# Row 23
$ws.Range("E23").Value = "  +0.36%  "

# Row 24
$ws.Range("E24").Value = "  -2.54%  "

# Row 25
$ws.Range("D25").Value = "'160.80"
$ws.Range("E25").Value = "  +0.87%  "

# Row 26
$ws.Range("E26").Value = "  +2.45%  "

# Row 27
$ws.Range("D27").Value = "'16.29"
$ws.Range("E27").Value = "  +0.35%  "

# Row 28
$ws.Range("E28").Value = "  +1.40%  "

# Row 29
$ws.Range("E29").Value = "  +0.05%  "

# Row 30
$ws.Range("E30").Value = "  +1.83%  "

# Row 31
$ws.Range("D31").Value = "'0.0520"
$ws.Range("E31").Value = "  +1.42%  "

# Row 32
$ws.Range("D32").Value = "'3.65"
$ws.Range("E32").Value = "  +0.06%  "

# Row 33
$ws.Range("E33").Value = "  +3.61%  "

# Row 34
$ws.Range("E34").Value = "  +1.61%  "

# Row 35
$ws.Range("D35").Value = "'1.406.03"
$ws.Range("E35").Value = "  +1.09%  "

# Row 36
$ws.Range("E36").Value = "  +1.31%  "

# Row 37
$ws.Range("D37").Value = "'0.0189"
$ws.Range("E37").Value = "  +2.40%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'2.37"
$ws.Range("E38").Value = "  +7.50%  "

# Row 39
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.04"
$ws.Range("E39").Value = "  -0.59%  "

# Row 40
$ws.Range("D40").Value = "'80.31"
$ws.Range("E40").Value = "  +2.65%  "

# Row 41
$ws.Range("E41").Value = "  +0.01%  "

# Row 42
$ws.Range("D42").Value = "'0.920"
$ws.Range("E42").Value = "  +1.53%  "

# Row 43
$ws.Range("E43").Value = "  -0.07%  "

# Row 44
$ws.Range("D44").Value = "'13.38"
$ws.Range("E44").Value = "  +9.96%  "

# Row 45
$ws.Range("E45").Value = "  +2.64%  "

# Row 46
$ws.Range("E46").Value = "  +3.32%  "

# Row 47
$ws.Range("E47").Value = "  -2.54%  "

# Row 48
$ws.Range("D48").Value = "'1.07"
$ws.Range("E48").Value = "  -0.19%  "

# Row 49
$ws.Range("D49").Value = "'106.95"
$ws.Range("E49").Value = "  -0.02%  "

# Row 50
$ws.Range("D50").Value = "'1.947.48"
$ws.Range("E50").Value = "  +0.05%  "

# Row 51
$ws.Range("E51").Value = "  +0.03%  "

Write-Host "Applied cryptos list update"